# EASY-1314: add SF_PLAY_MODE column to the instructions sheet.
#
# A new column is inserted immediately before the existing "AV_FILE"
# column (which sits at AV). The insert pushes AV:AZ -> AW:BA and the
# freshly created AV column picks up the formatting of its left
# neighbour (AU), matching Excel's normal "insert column" behaviour.
# The new column's header cell gets the value "SF_PLAY_MODE".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the width of the column that will end up to the left of the
# new one, so the inserted column can inherit it explicitly (Excel
# carries the left-hand column's width onto a freshly inserted column).
$leftWidth = $ws.Columns("AU:AU").ColumnWidth()

# Insert the new column at AV; everything from AV onward shifts right.
$ws.Columns("AV:AV").Insert() | Out-Null

# Give the inherited width to the new column explicitly.
$ws.Columns("AV:AV").ColumnWidth = $leftWidth

# Write the new header.
$ws.Range("AV1").Value = "SF_PLAY_MODE"

# The filter database named range covered A1:AY4; it must grow to keep
# covering the same logical columns now that one was inserted inside it.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$AZ`$4"
    }
}

# Leave the new header cell selected, as the last-edited cell.
$ws.Range("AV1").Select() | Out-Null
